$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) cells touched below to remain text,
# matching the source workbook where Price is stored as inlineStr text
# (values like "1.00" / "18.80" would otherwise be auto-coerced to numbers).
$priceCells = @("D10","D13","D14","D15","D16","D17","D19","D2","D20","D21","D23","D24","D25","D27","D28","D29","D3","D32","D33","D35","D36","D37","D4","D41","D44","D46","D47","D48","D49","D5","D50","D6","D7","D8","D9")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '51.918.67'
$ws.Range('D3').Value = '3.111.47'
$ws.Range('E3').Value = '  +4.20%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '389.63'
$ws.Range('E5').Value = '  +2.09%  '
$ws.Range('D6').Value = '104.10'
$ws.Range('E6').Value = '  +0.20%  '
$ws.Range('D7').Value = '0.546'
$ws.Range('E7').Value = '  +0.25%  '
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').Value = '0.593'
$ws.Range('E9').Value = '  +0.06%  '
$ws.Range('D10').Value = '37.30'
$ws.Range('E10').Value = '  +1.68%  '
$ws.Range('E11').Value = '  +0.45%  '
$ws.Range('E12').Value = '  +0.90%  '
$ws.Range('D13').Value = '3.587.77'
$ws.Range('E13').Value = '  +3.87%  '
$ws.Range('D14').Value = '18.80'
$ws.Range('E14').Value = '  +1.87%  '
$ws.Range('D15').Value = '7.89'
$ws.Range('E15').Value = '  +1.19%  '
$ws.Range('D16').Value = '3.099.92'
$ws.Range('E16').Value = '  +3.97%  '
$ws.Range('D17').Value = '0.987'
$ws.Range('E17').Value = '  -1.06%  '
$ws.Range('E18').Value = '  -4.19%  '
$ws.Range('D19').Value = '52.040.73'
$ws.Range('E19').Value = '  +1.08%  '
$ws.Range('D20').Value = '3.19'
$ws.Range('E20').Value = '  +2.86%  '
$ws.Range('D21').Value = '12.58'
$ws.Range('E21').Value = '  -0.21%  '
$ws.Range('E22').Value = '  +0.95%  '
$ws.Range('D23').Value = '70.52'
$ws.Range('E23').Value = '  +0.31%  '
$ws.Range('D24').Value = '270.08'
$ws.Range('E24').Value = '  +1.19%  '
$ws.Range('D25').Value = '3.15'
$ws.Range('E25').Value = '  -2.37%  '
$ws.Range('E26').Value = '  +5.31%  '
$ws.Range('D27').Value = '27.12'
$ws.Range('E27').Value = '  +3.89%  '
$ws.Range('D28').Value = '0.173'
$ws.Range('E28').Value = '  +2.53%  '
$ws.Range('D29').Value = '7.32'
$ws.Range('E29').Value = '  -0.33%  '
$ws.Range('E30').Value = '  +0.14%  '
$ws.Range('E31').Value = '  -0.54%  '
$ws.Range('D32').Value = '10.36'
$ws.Range('E32').Value = '  -0.42%  '
$ws.Range('D33').Value = '35.75'
$ws.Range('E33').Value = '  +2.66%  '
$ws.Range('E34').Value = '  +0.45%  '
$ws.Range('B35').Value = 'VeChain'
$ws.Range('C35').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D35').Value = '0.0453'
$ws.Range('E35').Value = '  +2.52%  '
$ws.Range('B36').Value = 'OKB'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D36').Value = '50.39'
$ws.Range('E36').Value = '  -1.99%  '
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  -0.12%  '
$ws.Range('E38').Value = '  +3.79%  '
$ws.Range('E39').Value = '  +8.94%  '
$ws.Range('E40').Value = '  +2.62%  '
$ws.Range('D41').Value = '17.04'
$ws.Range('E41').Value = '  +1.16%  '
$ws.Range('E42').Value = '  +0.20%  '
$ws.Range('E43').Value = '  -0.31%  '
$ws.Range('D44').Value = '127.49'
$ws.Range('E44').Value = '  +0.20%  '
$ws.Range('E45').Value = '  -1.13%  '
$ws.Range('D46').Value = '22.19'
$ws.Range('E46').Value = '  +3.64%  '
$ws.Range('B47').Value = 'WEMIXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D47').Value = '2.10'
$ws.Range('E47').Value = '  +2.99%  '
$ws.Range('B48').Value = 'ApeXProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D48').Value = '2.45'
$ws.Range('E48').Value = '  +3.54%  '
$ws.Range('D49').Value = '2.057.78'
$ws.Range('E49').Value = '  +1.43%  '
$ws.Range('D50').Value = '3.406.54'
$ws.Range('E50').Value = '  +3.84%  '
$ws.Range('E51').Value = '  +7.75%  '
